# Word COM-interop script implementing the commit diff for
# "Super Solutions for Shiny Architecture 1 of 5 Using Session Data.docx"
#
# Summary of textual changes:
# 1. Introduction paragraph: remove 'The official RStudio article
#    "[Communication between modules]" proposes ... alternative one.'
#    leaving just the sentence ending at "...modularized applications. "
# 2. "The Solution" paragraph: remove ' (check [this Joe Cheng post] if
#    you're not familiar with the idea and the usage of Shiny modules)'
#    leaving "...input, output, session. As each module has access to..."
# 3. Remove "Let's present the idea by modifying the [example presented
#    by RStudio]. " so the paragraph now starts with "In this rather
#    simple app..."
# 4. Remove the trailing "Article [Super Solutions...] comes from
#    [Appsilon...]." paragraph plus the empty paragraph that followed it.

$d = $word.ActiveDocument

function Find-End($needle) {
    # Executes a literal Find from the start of the document and returns
    # the End position of the first match.
    $r = $d.Content
    [void]$r.Find.Execute($needle, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    return $r.End
}

function Find-Start($needle) {
    # Executes a literal Find from the start of the document and returns
    # the Start position of the first match.
    $r = $d.Content
    [void]$r.Find.Execute($needle, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    return $r.Start
}

# --- Edit 1 -------------------------------------------------------------
$from = Find-End "multi-view, modularized applications."
$to = Find-End "our alternative one."
$r = $d.Range($from, $to)
$r.Text = " "

# --- Edit 2 ---------------------------------------------------------------
# Delete the " (check [hyperlink] if you're not familiar ... modules). "
# span (including the trailing space before "As"), then prepend ". " onto
# the surviving "As each module..." run so the new punctuation keeps that
# run's (non-italic) formatting rather than bleeding into the preceding
# italic "input, output, session" run.
$from = Find-Start " (check "
$to = Find-End "if you're not familiar with the idea and the usage of Shiny modules). "
$r = $d.Range($from, $to)
$r.Delete()

$target = $d.Range($from, $from + 1)
$target.Text = ". " + $target.Text

# --- Edit 3 -------------------------------------------------------------
$from = Find-Start "Let's present the idea by modifying the "
$to = Find-Start "In this rather simple "
$r = $d.Range($from, $to)
$r.Text = ""

# --- Edit 4 ---------------------------------------------------------------
$pCount = $d.Paragraphs.Count
$pSumUp = $d.Paragraphs.Item($pCount - 2)
$pEmpty = $d.Paragraphs.Item($pCount)
$r = $d.Range($pSumUp.Range.End, $pEmpty.Range.End)
$r.Delete()
